$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Locations")

$values = @(
  "Toronto",
  "Locations",
  "New York",
  "California",
  "London",
  "Paris",
  "Madrid",
  "Milan",
  "Barcelona",
  "Lisbon",
  "Rome",
  "Manchester",
  "Newcastle",
  "Sydney",
  "Perth",
  "Liverpool",
  "Berlin"
)

# Type the values in order down column A (A1:A17), matching the
# existing "black font" style used by the original two rows.
for ($i = 0; $i -lt $values.Length; $i++) {
  $row = $i + 1
  $cell = $ws.Cells.Item($row, 1)
  $cell.Value = $values[$i]
  $cell.Font.Color = 0
}

# Swap A1 and A2 so the "Locations" header ends up on top, above "Toronto".
$ws.Cells.Item(1, 1).Value = "Locations"
$ws.Cells.Item(2, 1).Value = "Toronto"

$ws.Range("A17").Select()
